$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1 (title slide): "Text Placeholder 6" -- update the tutorial title /
# venue line from the NOAA Global Systems Laboratory wording to the
# "Better Software for Reproducible Science tutorial @ SC23" wording, with
# the first sentence styled as its own run (explicit color / theme font).
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$shp1 = $slide1.Shapes.Item(6)
$tf1 = $shp1.TextFrame
$tr1 = $tf1.TextRange

$newTitleLine = "Better Software for Reproducible Science tutorial @ SC23"
$tr1.Text = $newTitleLine

# Remove the extra paragraph spacing that preceded this paragraph.
$tr1.ParagraphFormat.SpaceBefore = 0

# First run: "Better Software for Reproducible Science tutorial"
$run1 = $tr1.Characters(1, 49)
$run1.Font.Bold = $false
$run1.Font.Italic = $false
$run1.Font.Shadow = $false
$run1.Font.Color.RGB = 0x111111
$run1.Font.Name = "+mn-lt"

# Second run: the single space separator.
$run2 = $tr1.Characters(50, 1)
$run2.Font.Name = "+mn-lt"

# Third run: "@ SC23" keeps the plain default formatting.

# ---------------------------------------------------------------------------
# Slide 2 (License, Citation and Acknowledgements): update the citation
# paragraph text (new author list / venue / DOI) and switch the citation
# runs from the major theme font to the minor theme font.
# ---------------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$shp2 = $slide2.Shapes.Item(2)
$tf2 = $shp2.TextFrame
$tr2 = $tf2.TextRange
$citationPara = $tr2.Paragraphs(3)

$citationRun = $citationPara.Runs(2)
$citationRun.Text = "David E. Bernholdt, Patricia A. Grubel, David M. Rogers, and Gregory R. Watson, Better Software for Reproducible Science tutorial, in The International Conference for High-Performance Computing, Networking, Storage, and Analysis (SC23), Denver, Colorado, 2023. DOI:" + [char]0x00A0
$citationRun.Font.Name = "+mn-lt"

$doiRun = $citationPara.Runs(3)
$doiRun.Text = "10.6084/m9.figshare.24226105"
$doiRun.Font.Name = "+mn-lt"

$periodRun = $citationPara.Runs(4)
$periodRun.Font.Name = "+mn-lt"
